$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert three new rows above the current row 13 (the thin blank separator
# row) to make room for two new BOM lines (connectors) plus one blank
# spacer row, pushing the "Web ID / Access ID / short-link" block and the
# trailing blank rows down by three.
$ws.Rows("13:15").Insert()

$ltrMark = [char]0x200E

# --- Row 14 first (matches the authoring order the workbook was built in) ---
$ws.Range("A14").Value2 = 1
$ws.Range("B14").Value2 = [string]::Concat("0436500200", $ltrMark)
$ws.Range("C14").Value2 = "Connector"
$ws.Range("D14").Value2 = "CONN HEADER 2POS 3MM RT ANG TIN"
$ws.Range("F14").Value2 = "N/A"

# --- Row 13 next ---
$ws.Range("A13").Value2 = 1
$ws.Range("B13").Value2 = [string]::Concat("0436500400", $ltrMark)
$ws.Range("D13").Value2 = "CONN HEADER 4POS 3MM RT ANG TIN"
$ws.Range("C13").Value2 = "Connector"
$ws.Range("F13").Value2 = "N/A"

# Link column (H) filled in last, row 13 before row 14.
$ws.Range("H13").Value2 = "https://www.digikey.ca/product-detail/en/0436500400/WM1862-ND/268991/?itemSeq=276688085"
$ws.Range("H14").Value2 = "https://www.digikey.ca/product-detail/en/0436500200/WM1860-ND/268989/?itemSeq=276688088"

# The inserted rows copied formatting (incl. the Channels column E and the
# Quantity/Package columns) from the row above -- clear the cells that
# should not exist on the new BOM rows / the blank spacer row beneath them.
$ws.Range("E13").Clear()
$ws.Range("E14").Clear()
$ws.Range("A15").Clear()
$ws.Range("C15:F15").Clear()

# Widen the Description column slightly to fit the new connector text.
$ws.Columns("D").ColumnWidth = 31.67

# Move the selection to match where the author left off.
$ws.Range("B16").Select() | Out-Null
